$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Hoppers Crossing" exposure-site row (row 2) down to row 3
# so a new entry can be added above it. Re-read values with Value2 (Value's
# getter is unreliable in this runtime) and write them into row 3, preserving
# the plain/default formatting already used by that data row.
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2
# Also normalize the date text's day formatting while moving it (08/02/21 -> 8/2/21)
$ws.Range("C3").Value = "6.40am - 7.15am  8/2/21"

# Populate row 2 with the new Broadmeadows exposure site entry
$ws.Range("A2").Value = "Broadmeadows"
$ws.Range("B2").Value = "Broadmeadows Central  (West side of shopping centre, fresh fruit and meat section)  1099/1168 Pascoe Vale Rd  Broadmeadows VIC 3047"
$ws.Range("C2").Value = "12:15pm - 1:15pm  9/2/2021"
$ws.Range("D2").Value = "Case attended fresh fruit and meat section on the west side of the shopping centre"
